$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.224.54'
$ws.Range("E2").Value = '  -4.72%  '
$ws.Range("D3").Value = '2.572.01'
$ws.Range("E3").Value = '  -3.76%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = "'507.77"
$ws.Range("E5").Value = '  -4.53%  '
$ws.Range("D6").Value = "'144.96"
$ws.Range("E6").Value = '  -7.30%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = "'0.573"
$ws.Range("E8").Value = '  -2.41%  '
$ws.Range("D9").Value = '2.580.77'
$ws.Range("E9").Value = '  -4.07%  '
$ws.Range("E10").Value = '  -4.62%  '
$ws.Range("E11").Value = '  -5.23%  '
$ws.Range("E12").Value = '  -5.44%  '
$ws.Range("E13").Value = '  -0.87%  '
$ws.Range("D14").Value = '3.017.12'
$ws.Range("E14").Value = '  -4.07%  '
$ws.Range("D15").Value = '58.192.84'
$ws.Range("E15").Value = '  -4.84%  '
$ws.Range("D16").Value = "'21.07"
$ws.Range("E16").Value = '  -4.83%  '
$ws.Range("E17").Value = '  -4.38%  '
$ws.Range("D18").Value = '2.579.15'
$ws.Range("E18").Value = '  -3.97%  '
$ws.Range("D19").Value = "'4.54"
$ws.Range("E19").Value = '  -5.49%  '
$ws.Range("D20").Value = "'342.63"
$ws.Range("E20").Value = '  -3.78%  '
$ws.Range("D21").Value = "'10.30"
$ws.Range("E21").Value = '  -4.16%  '
$ws.Range("D22").Value = "'6.05"
$ws.Range("E22").Value = '  -4.90%  '
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = '  -0.29%  '
$ws.Range("D24").Value = "'60.47"
$ws.Range("E24").Value = '  -2.24%  '
$ws.Range("E25").Value = '  -3.59%  '
$ws.Range("E26").Value = '  -0.25%  '
$ws.Range("D27").Value = '2.680.74'
$ws.Range("E27").Value = '  -4.15%  '
$ws.Range("E28").Value = '  -5.86%  '
$ws.Range("D29").Value = '0.0₃0812'
$ws.Range("E29").Value = '  -6.32%  '
$ws.Range("E30").Value = '  -5.44%  '
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("D32").Value = "'6.10"
$ws.Range("E32").Value = '  -1.24%  '
$ws.Range("D33").Value = "'18.79"
$ws.Range("E33").Value = '  -4.20%  '
$ws.Range("D34").Value = "'148.93"
$ws.Range("E34").Value = '  -0.65%  '
$ws.Range("D35").Value = "'1.54"
$ws.Range("E35").Value = '  -5.82%  '
$ws.Range("E36").Value = '  +6.24%  '
$ws.Range("E37").Value = '  -4.72%  '
$ws.Range("E38").Value = '  -6.48%  '
$ws.Range("E39").Value = '  -6.91%  '
$ws.Range("E40").Value = '  -2.35%  '
$ws.Range("D41").Value = "'290.03"
$ws.Range("E41").Value = '  -5.46%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").Value = "'3.58"
$ws.Range("E42").Value = '  -5.69%  '
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").Value = "'1.40"
$ws.Range("E43").Value = '  -6.68%  '
$ws.Range("D44").Value = "'0.0990"
$ws.Range("E44").Value = '  -2.96%  '
$ws.Range("D45").Value = "'0.995"
$ws.Range("E45").Value = '  -0.22%  '
$ws.Range("E46").Value = '  -6.79%  '
$ws.Range("E47").Value = '  -5.68%  '
$ws.Range("D48").Value = "'19.17"
$ws.Range("E48").Value = '  -7.44%  '
$ws.Range("E49").Value = '  -1.00%  '
$ws.Range("D50").Value = "'0.0229"
$ws.Range("E50").Value = '  -5.49%  '
$ws.Range("D51").Value = "'4.57"
$ws.Range("E51").Value = '  -7.21%  '
